# The data row for E=7 (N19-0676, rep 2, plot 39) was accidentally
# duplicated (rows 107 and 108 were identical). Remove the duplicate
# row 108 so every subsequent record shifts up by one row and the
# lsmeans/export data lines back up correctly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(108).Delete()
